$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "SamplesTab" row (row 3) query text in column B had two columns
# dropped from its SELECT list: smp.sample_tumor_status AS Tumor and
# smp.sample_type AS "Analyte Type".
$newSamplesQuery = @"
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
  s.phs_accession = 'phs002529' AND d.primary_diagnosis = 'Yolk sac tumor'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
"@

$ws.Range("B3").Value = $newSamplesQuery

# Move the selection / active cell from C11 to B3, and scroll the sheet so
# B3 becomes the top-left visible cell.
$ws.Range("B3").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 2
